$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 'b'
$ws.Range("J2").Value = 'Acknowledge (Backchannel)'
$ws.Range("I7").Value = 'b'
$ws.Range("J7").Value = 'Acknowledge (Backchannel)'
$ws.Range("I12").Value = 'aa'
$ws.Range("J12").Value = 'Agree/Accept'
$ws.Range("I14").Value = 'aa'
$ws.Range("J14").Value = 'Agree/Accept'
$ws.Range("I28").Value = 'aa'
$ws.Range("J28").Value = 'Agree/Accept'
$ws.Range("I37").Value = 'sv'
$ws.Range("J37").Value = 'Statement-opinion'
$ws.Range("I39").Value = 'aa'
$ws.Range("J39").Value = 'Agree/Accept'
$ws.Range("I45").Value = 'sd'
$ws.Range("J45").Value = 'Statement-non-opinion'
$ws.Range("I59").Value = '%'
$ws.Range("J59").Value = 'Uninterpretable'
$ws.Range("I62").Value = 'aa'
$ws.Range("J62").Value = 'Agree/Accept'
$ws.Range("I63").Value = 'sv'
$ws.Range("J63").Value = 'Statement-opinion'
$ws.Range("I70").Value = 'aa'
$ws.Range("J70").Value = 'Agree/Accept'
$ws.Range("I73").Value = 'sv'
$ws.Range("J73").Value = 'Statement-opinion'
$ws.Range("I74").Value = 'sd'
$ws.Range("J74").Value = 'Statement-non-opinion'
$ws.Range("I75").Value = 'b'
$ws.Range("J75").Value = 'Acknowledge (Backchannel)'
$ws.Range("I77").Value = 'sv'
$ws.Range("J77").Value = 'Statement-opinion'
$ws.Range("I95").Value = 'sd'
$ws.Range("J95").Value = 'Statement-non-opinion'
$ws.Range("I123").Value = 'sd'
$ws.Range("J123").Value = 'Statement-non-opinion'
$ws.Range("I128").Value = 'sd'
$ws.Range("J128").Value = 'Statement-non-opinion'
$ws.Range("I136").Value = 'b'
$ws.Range("J136").Value = 'Acknowledge (Backchannel)'
$ws.Range("I137").Value = 'aa'
$ws.Range("J137").Value = 'Agree/Accept'
$ws.Range("I156").Value = 'sd'
$ws.Range("J156").Value = 'Statement-non-opinion'
$ws.Range("I158").Value = 'sd'
$ws.Range("J158").Value = 'Statement-non-opinion'
$ws.Range("I169").Value = 'sd'
$ws.Range("J169").Value = 'Statement-non-opinion'
$ws.Range("I172").Value = '%'
$ws.Range("J172").Value = 'Uninterpretable'
$ws.Range("I173").Value = 'sv'
$ws.Range("J173").Value = 'Statement-opinion'
$ws.Range("I176").Value = 'sd'
$ws.Range("J176").Value = 'Statement-non-opinion'
$ws.Range("I186").Value = 'sv'
$ws.Range("J186").Value = 'Statement-opinion'
$ws.Range("I190").Value = 'sv'
$ws.Range("J190").Value = 'Statement-opinion'
$ws.Range("I196").Value = 'b'
$ws.Range("J196").Value = 'Acknowledge (Backchannel)'
$ws.Range("I222").Value = 'sv'
$ws.Range("J222").Value = 'Statement-opinion'
$ws.Range("I227").Value = '%'
$ws.Range("J227").Value = 'Uninterpretable'
$ws.Range("I230").Value = 'sd'
$ws.Range("J230").Value = 'Statement-non-opinion'
$ws.Range("I236").Value = 'aa'
$ws.Range("J236").Value = 'Agree/Accept'
$ws.Range("I249").Value = 'sd'
$ws.Range("J249").Value = 'Statement-non-opinion'
$ws.Range("I251").Value = 'sv'
$ws.Range("J251").Value = 'Statement-opinion'
$ws.Range("I261").Value = '%'
$ws.Range("J261").Value = 'Uninterpretable'
$ws.Range("I264").Value = 'sv'
$ws.Range("J264").Value = 'Statement-opinion'
$ws.Range("I266").Value = 'sd'
$ws.Range("J266").Value = 'Statement-non-opinion'
$ws.Range("I289").Value = 'sv'
$ws.Range("J289").Value = 'Statement-opinion'
$ws.Range("I296").Value = 'sd'
$ws.Range("J296").Value = 'Statement-non-opinion'
$ws.Range("I305").Value = '%'
$ws.Range("J305").Value = 'Uninterpretable'
$ws.Range("I306").Value = '%'
$ws.Range("J306").Value = 'Uninterpretable'
$ws.Range("I319").Value = '%'
$ws.Range("J319").Value = 'Uninterpretable'
$ws.Range("I321").Value = 'sv'
$ws.Range("J321").Value = 'Statement-opinion'
$ws.Range("I330").Value = 'sd'
$ws.Range("J330").Value = 'Statement-non-opinion'
$ws.Range("I340").Value = 'sv'
$ws.Range("J340").Value = 'Statement-opinion'
$ws.Range("I346").Value = '%'
$ws.Range("J346").Value = 'Uninterpretable'
$ws.Range("I367").Value = '%'
$ws.Range("J367").Value = 'Uninterpretable'
$ws.Range("I370").Value = 'sd'
$ws.Range("J370").Value = 'Statement-non-opinion'
$ws.Range("I375").Value = 'sv'
$ws.Range("J375").Value = 'Statement-opinion'
$ws.Range("I376").Value = 'sv'
$ws.Range("J376").Value = 'Statement-opinion'
$ws.Range("I395").Value = 'b'
$ws.Range("J395").Value = 'Acknowledge (Backchannel)'
$ws.Range("I427").Value = '%'
$ws.Range("J427").Value = 'Uninterpretable'
$ws.Range("I432").Value = 'aa'
$ws.Range("J432").Value = 'Agree/Accept'
$ws.Range("I433").Value = 'aa'
$ws.Range("J433").Value = 'Agree/Accept'
